# Insert 3 new daily-price rows ("Especial"/"Primera"/"Segunda" for
# Cultivar IV Región, date 2021-12-16 / serial 44546) at the top of the
# data block (row 78), pushing all the existing rows down by three.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("78:80").Insert()

# --- Columns shared by all three new rows ---
$ws.Range("A78:A80").Value = 3
$ws.Range("B78:B80").Value = "Femacal de La Calera"
$ws.Range("C78:C80").Value = "Coquimbo"
$ws.Range("D78:D80").Value = 44546
$ws.Range("E78:E80").Value = 5
$ws.Range("F78:F80").Value = "Fruta"
$ws.Range("G78:G80").Value = 100107
$ws.Range("H78:H80").Value = "Otros"
$ws.Range("I78:I80").Value = 100107002
$ws.Range("J78:J80").Value = "Chirimoya"
$ws.Range("K78:K80").Value = "Cultivar IV Región"
$ws.Range("Q78:Q80").Value = "`$/bandeja 10 kilos"
$ws.Range("R78:R80").Value = "Provincia de Limarí"
$ws.Range("T78:T80").Value = 10

# --- Row 78: Especial ---
$ws.Range("L78").Value = "Especial"
$ws.Range("M78").Value = 65
$ws.Range("N78").Value = 24000
$ws.Range("O78").Value = 24000
$ws.Range("P78").Value = 24000
$ws.Range("S78").Value = 2400

# --- Row 79: Primera ---
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 67
$ws.Range("N79").Value = 22000
$ws.Range("O79").Value = 22000
$ws.Range("P79").Value = 22000
$ws.Range("S79").Value = 2200

# --- Row 80: Segunda ---
$ws.Range("L80").Value = "Segunda"
$ws.Range("M80").Value = 60
$ws.Range("N80").Value = 20000
$ws.Range("O80").Value = 20000
$ws.Range("P80").Value = 20000
$ws.Range("S80").Value = 2000
